$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.8029023333333333
$ws.Range("H2").Value = 2.408707
$ws.Range("I2").Value = 0.06206726394886004
$ws.Range("J2").Value = 0.06206726394886004
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 97.419871
$ws.Range("N2").Value = 292.259613
$ws.Range("O2").Value = 0.3451284562056485
$ws.Range("P2").Value = 0.3451284562056485
$ws.Range("Q2").Value = 78.21864173893233
$ws.Range("R2").Value = 703.9677756503909
$ws.Range("S2").Value = 0.02142117898757857
$ws.Range("T2").Value = 0.02142117898757857

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.8029023333333333
$ws.Range("H3").Value = 2.408707
$ws.Range("I3").Value = 0.06206726394886004
$ws.Range("J3").Value = 0.06206726394886004
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 34.912838
$ws.Range("N3").Value = 104.738514
$ws.Range("O3").Value = 0.123685381195977
$ws.Range("P3").Value = 0.123685381195977
$ws.Range("Q3").Value = 28.03159909348867
$ws.Range("R3").Value = 252.284391841398
$ws.Range("S3").Value = 0.007676813201306074
$ws.Range("T3").Value = 0.007676813201306073

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.8029023333333333
$ws.Range("H4").Value = 2.408707
$ws.Range("I4").Value = 0.06206726394886004
$ws.Range("J4").Value = 0.06206726394886004
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 50.427193
$ws.Range("N4").Value = 151.281579
$ws.Range("O4").Value = 0.1786479400170247
$ws.Range("P4").Value = 0.1786479400170247
$ws.Range("Q4").Value = 40.48811092315033
$ws.Range("R4").Value = 364.392998308353
$ws.Range("S4").Value = 0.01108818884695679
$ws.Range("T4").Value = 0.01108818884695679

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.8029023333333333
$ws.Range("H5").Value = 2.408707
$ws.Range("I5").Value = 0.06206726394886004
$ws.Range("J5").Value = 0.06206726394886004
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 63.884013
$ws.Range("N5").Value = 191.652039
$ws.Range("O5").Value = 0.2263212890408322
$ws.Range("P5").Value = 0.2263212890408321
$ws.Range("Q5").Value = 51.292623100397
$ws.Range("R5").Value = 461.633607903573
$ws.Range("S5").Value = 0.01404714318414358
$ws.Range("T5").Value = 0.01404714318414357

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.8029023333333333
$ws.Range("H6").Value = 2.408707
$ws.Range("I6").Value = 0.06206726394886004
$ws.Range("J6").Value = 0.06206726394886004
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 19.609342
$ws.Range("N6").Value = 58.82802600000001
$ws.Range("O6").Value = 0.06946983056124746
$ws.Range("P6").Value = 0.06946983056124745
$ws.Range("Q6").Value = 15.74438644693133
$ws.Range("R6").Value = 141.699478022382
$ws.Range("S6").Value = 0.004311802309927531
$ws.Range("T6").Value = 0.00431180230992753

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.8029023333333333
$ws.Range("H7").Value = 2.408707
$ws.Range("I7").Value = 0.06206726394886004
$ws.Range("J7").Value = 0.06206726394886004
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 16.01808066666667
$ws.Range("N7").Value = 48.054242
$ws.Range("O7").Value = 0.05674710297927013
$ws.Range("P7").Value = 0.05674710297927013
$ws.Range("Q7").Value = 12.86095434278822
$ws.Range("R7").Value = 115.748589085094
$ws.Range("S7").Value = 0.003522137418947502
$ws.Range("T7").Value = 0.003522137418947502

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.159929333333333
$ws.Range("H8").Value = 6.479788
$ws.Range("I8").Value = 0.1669703754456877
$ws.Range("J8").Value = 0.1669703754456877
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 97.419871
$ws.Range("N8").Value = 292.259613
$ws.Range("O8").Value = 0.3451284562056485
$ws.Range("P8").Value = 0.3451284562056485
$ws.Range("Q8").Value = 210.4200370224493
$ws.Range("R8").Value = 1893.780333202044
$ws.Range("S8").Value = 0.05762622790964771
$ws.Range("T8").Value = 0.05762622790964769

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.159929333333333
$ws.Range("H9").Value = 6.479788
$ws.Range("I9").Value = 0.1669703754456877
$ws.Range("J9").Value = 0.1669703754456877
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 34.912838
$ws.Range("N9").Value = 104.738514
$ws.Range("O9").Value = 0.123685381195977
$ws.Range("P9").Value = 0.123685381195977
$ws.Range("Q9").Value = 75.40926290611466
$ws.Range("R9").Value = 678.6833661550321
$ws.Range("S9").Value = 0.02065179453543527
$ws.Range("T9").Value = 0.02065179453543527

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.159929333333333
$ws.Range("H10").Value = 6.479788
$ws.Range("I10").Value = 0.1669703754456877
$ws.Range("J10").Value = 0.1669703754456877
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 50.427193
$ws.Range("N10").Value = 151.281579
$ws.Range("O10").Value = 0.1786479400170247
$ws.Range("P10").Value = 0.1786479400170247
$ws.Range("Q10").Value = 108.9191733583613
$ws.Range("R10").Value = 980.2725602252519
$ws.Range("S10").Value = 0.0298289136172413
$ws.Range("T10").Value = 0.0298289136172413

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.159929333333333
$ws.Range("H11").Value = 6.479788
$ws.Range("I11").Value = 0.1669703754456877
$ws.Range("J11").Value = 0.1669703754456877
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 63.884013
$ws.Range("N11").Value = 191.652039
$ws.Range("O11").Value = 0.2263212890408322
$ws.Range("P11").Value = 0.2263212890408321
$ws.Range("Q11").Value = 137.984953609748
$ws.Range("R11").Value = 1241.864582487732
$ws.Range("S11").Value = 0.03778895060249974
$ws.Range("T11").Value = 0.03778895060249973

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.159929333333333
$ws.Range("H12").Value = 6.479788
$ws.Range("I12").Value = 0.1669703754456877
$ws.Range("J12").Value = 0.1669703754456877
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 19.609342
$ws.Range("N12").Value = 58.82802600000001
$ws.Range("O12").Value = 0.06946983056124746
$ws.Range("P12").Value = 0.06946983056124745
$ws.Range("Q12").Value = 42.35479299316534
$ws.Range("R12").Value = 381.1931369384881
$ws.Range("S12").Value = 0.0115994036909598
$ws.Range("T12").Value = 0.01159940369095979

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.159929333333333
$ws.Range("H13").Value = 6.479788
$ws.Range("I13").Value = 0.1669703754456877
$ws.Range("J13").Value = 0.1669703754456877
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 16.01808066666667
$ws.Range("N13").Value = 48.054242
$ws.Range("O13").Value = 0.05674710297927013
$ws.Range("P13").Value = 0.05674710297927013
$ws.Range("Q13").Value = 34.59792229563289
$ws.Range("R13").Value = 311.381300660696
$ws.Range("S13").Value = 0.009475085089903834
$ws.Range("T13").Value = 0.009475085089903834

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 9.973171666666666
$ws.Range("H14").Value = 29.919515
$ws.Range("I14").Value = 0.7709623606054523
$ws.Range("J14").Value = 0.7709623606054523
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 97.419871
$ws.Range("N14").Value = 292.259613
$ws.Range("O14").Value = 0.3451284562056485
$ws.Range("P14").Value = 0.3451284562056485
$ws.Range("Q14").Value = 971.5850972275216
$ws.Range("R14").Value = 8744.265875047695
$ws.Range("S14").Value = 0.2660810493084222
$ws.Range("T14").Value = 0.2660810493084222

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 9.973171666666666
$ws.Range("H15").Value = 29.919515
$ws.Range("I15").Value = 0.7709623606054523
$ws.Range("J15").Value = 0.7709623606054523
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 34.912838
$ws.Range("N15").Value = 104.738514
$ws.Range("O15").Value = 0.123685381195977
$ws.Range("P15").Value = 0.123685381195977
$ws.Range("Q15").Value = 348.1917267445233
$ws.Range("R15").Value = 3133.72554070071
$ws.Range("S15").Value = 0.09535677345923563
$ws.Range("T15").Value = 0.09535677345923561

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 9.973171666666666
$ws.Range("H16").Value = 29.919515
$ws.Range("I16").Value = 0.7709623606054523
$ws.Range("J16").Value = 0.7709623606054523
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 50.427193
$ws.Range("N16").Value = 151.281579
$ws.Range("O16").Value = 0.1786479400170247
$ws.Range("P16").Value = 0.1786479400170247
$ws.Range("Q16").Value = 502.9190524571316
$ws.Range("R16").Value = 4526.271472114185
$ws.Range("S16").Value = 0.1377308375528266
$ws.Range("T16").Value = 0.1377308375528266

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 9.973171666666666
$ws.Range("H17").Value = 29.919515
$ws.Range("I17").Value = 0.7709623606054523
$ws.Range("J17").Value = 0.7709623606054523
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 63.884013
$ws.Range("N17").Value = 191.652039
$ws.Range("O17").Value = 0.2263212890408322
$ws.Range("P17").Value = 0.2263212890408321
$ws.Range("Q17").Value = 637.1262284045649
$ws.Range("R17").Value = 5734.136055641085
$ws.Range("S17").Value = 0.1744851952541888
$ws.Range("T17").Value = 0.1744851952541888

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 9.973171666666666
$ws.Range("H18").Value = 29.919515
$ws.Range("I18").Value = 0.7709623606054523
$ws.Range("J18").Value = 0.7709623606054523
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 19.609342
$ws.Range("N18").Value = 58.82802600000001
$ws.Range("O18").Value = 0.06946983056124746
$ws.Range("P18").Value = 0.06946983056124745
$ws.Range("Q18").Value = 195.5673340363766
$ws.Range("R18").Value = 1760.10600632739
$ws.Range("S18").Value = 0.05355862456036014
$ws.Range("T18").Value = 0.05355862456036013

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 9.973171666666666
$ws.Range("H19").Value = 29.919515
$ws.Range("I19").Value = 0.7709623606054523
$ws.Range("J19").Value = 0.7709623606054523
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 16.01808066666667
$ws.Range("N19").Value = 48.054242
$ws.Range("O19").Value = 0.05674710297927013
$ws.Range("P19").Value = 0.05674710297927013
$ws.Range("Q19").Value = 159.7510682591811
$ws.Range("R19").Value = 1437.75961433263
$ws.Range("S19").Value = 0.0437498804704188
$ws.Range("T19").Value = 0.0437498804704188
